# Update contribution values to reflect new supply side/IRA assumptions.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    10 = @{ 'I' = 0.0969;  'J' = -0.045;   'K' = 0.2751;  'L' = 0.2384;  'M' = 0.5026;  'N' = 0.0985;   'O' = 0.087 }
    11 = @{ 'P' = 0.0027;  'Q' = -0.0924;  'R' = -0.0876; 'S' = -0.1505; 'T' = -0.0734; 'U' = -0.0804;  'V' = -0.0891; 'W' = -0.0788; 'X' = -1.0948 }
    30 = @{ 'I' = -2.2972; 'J' = -0.5605;  'K' = 0.0439;  'L' = -0.2074; 'M' = 0.7225;  'N' = 0.0501;   'O' = -0.3236 }
    31 = @{ 'P' = -0.6367; 'Q' = -0.309;   'R' = -0.0824; 'S' = -0.6647; 'T' = -0.7811; 'U' = -0.5699;  'V' = -0.5089; 'W' = -0.124;  'X' = -72.4415 }
    66 = @{ 'I' = -0.0098; 'J' = 0.0002;   'K' = -0.0199; 'L' = -0.019;  'M' = 0.1938;  'N' = -0.1811;  'O' = -0.0691 }
    67 = @{ 'P' = 0.1229;  'Q' = 0.1457;   'R' = 0.1359;  'S' = 0.1334;  'T' = 0.146;   'U' = 0.0903;   'V' = -0.0537; 'W' = 0.0354;  'X' = -0.6676 }
    86 = @{ 'I' = -0.0098; 'J' = 0.0002;   'K' = -0.0199; 'L' = -0.019;  'M' = 0.1938;  'N' = -0.1811;  'O' = -0.0691 }
    87 = @{ 'P' = 0.0371;  'Q' = 0.054;    'R' = 0.0367;  'S' = 0.0215;  'T' = 0.0384;  'U' = -0.0057;  'V' = -0.1384; 'W' = -0.0331; 'X' = -1.7944 }
}

foreach ($rowNum in $updates.Keys) {
    $rowUpdates = $updates[$rowNum]
    foreach ($col in $rowUpdates.Keys) {
        $cellRef = "$col$rowNum"
        $ws.Range($cellRef).Value = $rowUpdates[$col]
    }
}
